$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so values such as "1.003" or
# "0.06712" are not silently converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.628.79'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '1.919.52'
$ws.Range("E3").Value = '  -1.92%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '239.07'
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("D7").Value = '0.4773'
$ws.Range("E7").Value = '  -2.10%  '
$ws.Range("D8").Value = '0.2879'
$ws.Range("E8").Value = '  -2.74%  '
$ws.Range("D9").Value = '0.06712'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").Value = '18.80'
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("D11").Value = '104.25'
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.921.63'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07726'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("D15").Value = '0.6857'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").Value = '266.14'
$ws.Range("E16").Value = '  -6.42%  '
$ws.Range("D17").Value = '30.646.97'
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '1.003'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '0.000007520'
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("D20").Value = '12.78'
$ws.Range("E20").Value = '  -3.36%  '
$ws.Range("D21").Value = '5.461'
$ws.Range("E21").Value = '  -0.91%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '6.352'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '9.693'
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '163.70'
$ws.Range("E25").Value = '  -3.92%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '19.11'
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.099'
$ws.Range("E27").Value = '  -5.17%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1023'
$ws.Range("E28").Value = '  -3.25%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.393'
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.523'
$ws.Range("E30").Value = '  -3.79%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.446'
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.262'
$ws.Range("E32").Value = '  -4.48%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.04765'
$ws.Range("E33").Value = '  -3.43%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.7374'
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.122'
$ws.Range("E35").Value = '  -4.26%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01958'
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.643'
$ws.Range("E39").Value = '  -2.04%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.359'
$ws.Range("E40").Value = '  -3.08%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '75.35'
$ws.Range("E41").Value = '  -2.72%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '2.005'
$ws.Range("E42").Value = '  -5.19%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8667'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '106.97'
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4321'
$ws.Range("E45").Value = '  -3.36%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.004'
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.586'
$ws.Range("E47").Value = '  -6.56%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '955.32'
$ws.Range("E48").Value = '  -4.23%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1206'
$ws.Range("E49").Value = '  -4.35%  '
$ws.Range("D50").Value = '35.24'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '8.966'
$ws.Range("E51").Value = '  -3.89%  '
